# Generate wordclouds for NER
# Adds two new named-entity rows (Evidence / Moon Landing) to the "2017"
# sheet, mirroring the existing "Lunar Roving Vehicle" row, including the
# Wikidata hyperlink in column D styled like the existing hyperlink cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: Evidence / Q1237694 -------------------------------------------------
$ws.Range("A3").Value2 = "Evidence"
$ws.Range("C3").Value2 = 1237694
$ws.Range("D3").Value2 = "https://www.wikidata.org/wiki/Q1237694"
$ws.Range("E3").Value2 = "Evidence"
$ws.Range("F3").Value2 = 1

# --- Row 4: Moon Landing / Q17175022 --------------------------------------------
$ws.Range("A4").Value2 = "Moon Landing"
$ws.Range("C4").Value2 = 17175022
$ws.Range("D4").Value2 = "https://www.wikidata.org/wiki/Q17175022"
$ws.Range("E4").Value2 = "Moon Landings"
$ws.Range("F4").Value2 = 1

# --- Hyperlinks for the new Wikidata_url cells, matching existing D2 style -----
$ws.Hyperlinks.Add($ws.Range("D3"), "https://www.wikidata.org/wiki/Q1237694")
$ws.Hyperlinks.Add($ws.Range("D4"), "https://www.wikidata.org/wiki/Q17175022")

$ws.Range("D3").Style = $ws.Range("D2").Style
$ws.Range("D4").Style = $ws.Range("D2").Style

$ws.Range("D3").Value2 = "https://www.wikidata.org/wiki/Q1237694"
$ws.Range("D4").Value2 = "https://www.wikidata.org/wiki/Q17175022"
